$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.86 = 6717.87 pesos`n✅ 6717.87 pesos = 1.85 = 953.46 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate cells on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 536.777
$ws2.Range("O10").Value = 3606
$ws2.Range("N12").Value = 3630
$ws2.Range("O12").Value = 515.2
